$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3974.889
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 3974.889
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 11924.667
$ws.Range("N70").Value = -12464.667
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 3974.889
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 3974.889
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 11924.667
$ws.Range("N73").Value = -13796.667
$ws.Range("M73").ClearContents()
$ws.Range("H138").Value = 1983.79
$ws.Range("I138").Value = 718.78125
$ws.Range("J138").Value = 2579.0881
$ws.Range("K138").Value = 2156.34375
$ws.Range("L138").Value = 7737.2643
$ws.Range("M138").Value = 2983.65625
$ws.Range("N138").Value = -18017.2643
$ws.Range("H141").Value = 5733.864
$ws.Range("I141").Value = 5844.75
$ws.Range("J141").Value = 4625
$ws.Range("K141").Value = 17534.25
$ws.Range("L141").Value = 13875
$ws.Range("M141").Value = -12354.25
$ws.Range("N141").Value = -24235

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 799.2692
$ws.Range("I2").Value = 668.45
$ws.Range("J2").Value = 1235.3334
$ws.Range("K2").Value = 668.45
$ws.Range("L2").Value = 1235.3334
$ws.Range("M2").Value = -555.45
$ws.Range("N2").Value = -1461.3334
$ws.Range("H32").Value = 6578
$ws.Range("I32").Value = 6083.795
$ws.Range("J32").Value = 8060.615
$ws.Range("K32").Value = 6083.795
$ws.Range("L32").Value = 8060.615
$ws.Range("M32").Value = -5796.795
$ws.Range("N32").Value = -8634.615
$ws.Range("H49").Value = 18333
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 18333
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 18333
$ws.Range("N49").Value = -18853
$ws.Range("H76").Value = 38700
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 38700
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 38700
$ws.Range("N76").Value = -39376
$ws.Range("H79").Value = 38700
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 38700
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 38700
$ws.Range("N79").Value = -41040
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H116").Value = 799.2692
$ws.Range("I116").Value = 668.45
$ws.Range("J116").Value = 1235.3334
$ws.Range("K116").Value = 668.45
$ws.Range("L116").Value = 1235.3334
$ws.Range("M116").Value = 1625.55
$ws.Range("N116").Value = -5823.3334
$ws.Range("H132").Value = 2403.5881
$ws.Range("I132").Value = 1231.6522
$ws.Range("J132").Value = 4854
$ws.Range("K132").Value = 3694.9566
$ws.Range("L132").Value = 14562
$ws.Range("M132").Value = -1164.9566
$ws.Range("N132").Value = -19622

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 799.2692
$ws.Range("I3").Value = 668.45
$ws.Range("J3").Value = 1235.3334
$ws.Range("K3").Value = 668.45
$ws.Range("L3").Value = 1235.3334
$ws.Range("M3").Value = -554.45
$ws.Range("N3").Value = -1463.3334
$ws.Range("H54").Value = 22999
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 22999
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 22999
$ws.Range("N54").Value = -23967

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1752.6329
$ws.Range("I58").Value = 1565.8308
$ws.Range("J58").Value = 2619.9285
$ws.Range("K58").Value = 1565.8308
$ws.Range("L58").Value = 2619.9285
$ws.Range("M58").Value = -1362.8308
$ws.Range("N58").Value = -3025.9285
$ws.Range("H62").Value = 7302
$ws.Range("I62").Value = 10000
$ws.Range("J62").Value = 5953
$ws.Range("K62").Value = 10000
$ws.Range("L62").Value = 5953
$ws.Range("M62").Value = -9376
$ws.Range("N62").Value = -7201
$ws.Range("H65").Value = 7302
$ws.Range("I65").Value = 10000
$ws.Range("J65").Value = 5953
$ws.Range("K65").Value = 50000
$ws.Range("L65").Value = 29765
$ws.Range("M65").Value = -46880
$ws.Range("N65").Value = -36005
$ws.Range("H134").Value = 5968.1665
$ws.Range("I134").Value = 9186.916999999999
$ws.Range("J134").Value = 2749.4167
$ws.Range("K134").Value = 27560.751
$ws.Range("L134").Value = 8248.250100000001
$ws.Range("M134").Value = -25025.751
$ws.Range("N134").Value = -13318.2501
$ws.Range("H136").Value = 1752.6329
$ws.Range("I136").Value = 1565.8308
$ws.Range("J136").Value = 2619.9285
$ws.Range("K136").Value = 4697.4924
$ws.Range("L136").Value = 7859.7855
$ws.Range("M136").Value = -2147.4924
$ws.Range("N136").Value = -12959.7855
$ws.Range("H138").Value = 42834
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 42834
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 42834
$ws.Range("N138").Value = -53114
$ws.Range("H140").Value = 73781.875
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 73781.875
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 73781.875
$ws.Range("N140").Value = -84141.875
$ws.Range("H141").Value = 31357.143
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 31357.143
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 31357.143
$ws.Range("N141").Value = -41717.143

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 854.6667
$ws.Range("I113").Value = 596.5
$ws.Range("J113").Value = 983.75
$ws.Range("K113").Value = 1789.5
$ws.Range("L113").Value = 2951.25
$ws.Range("M113").Value = 380.5
$ws.Range("N113").Value = -7291.25
$ws.Range("H131").Value = 8621593
$ws.Range("I131").Value = 100000280
$ws.Range("J131").Value = 961.8679
$ws.Range("K131").Value = 300000840
$ws.Range("L131").Value = 2885.6037
$ws.Range("M131").Value = -299995800
$ws.Range("N131").Value = -12965.6037

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 14299.75
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 14299.75
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 14299.75
$ws.Range("N39").Value = -15363.75
$ws.Range("H132").Value = 3364.9167
$ws.Range("I132").Value = 2208
$ws.Range("J132").Value = 4191.2856
$ws.Range("K132").Value = 6624
$ws.Range("L132").Value = 12573.8568
$ws.Range("M132").Value = -4094
$ws.Range("N132").Value = -17633.8568
$ws.Range("H140").Value = 38622.332
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 38622.332
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 38622.332
$ws.Range("N140").Value = -48982.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8738.200000000001
$ws.Range("I40").Value = 9044.875
$ws.Range("J40").Value = 8533.75
$ws.Range("K40").Value = 9044.875
$ws.Range("L40").Value = 8533.75
$ws.Range("M40").Value = -8908.875
$ws.Range("N40").Value = -8805.75
$ws.Range("H46").Value = 2703.4443
$ws.Range("I46").Value = 2819.75
$ws.Range("J46").Value = 2610.4
$ws.Range("K46").Value = 2819.75
$ws.Range("L46").Value = 2610.4
$ws.Range("M46").Value = -2631.75
$ws.Range("N46").Value = -2986.4
$ws.Range("H61").Value = 1721.5
$ws.Range("I61").Value = 1704.1666
$ws.Range("J61").Value = 1747.5
$ws.Range("K61").Value = 1704.1666
$ws.Range("L61").Value = 1747.5
$ws.Range("M61").Value = -1502.1666
$ws.Range("N61").Value = -2151.5
$ws.Range("H74").Value = 43202.832
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 43202.832
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 43202.832
$ws.Range("N74").Value = -45198.832
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 43202.832
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 43202.832
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 129608.496
$ws.Range("N77").Value = -139592.496
$ws.Range("M77").ClearContents()
$ws.Range("H113").Value = 1721.5
$ws.Range("I113").Value = 1704.1666
$ws.Range("J113").Value = 1747.5
$ws.Range("K113").Value = 1704.1666
$ws.Range("L113").Value = 1747.5
$ws.Range("M113").Value = 465.8334
$ws.Range("N113").Value = -6087.5
$ws.Range("H122").Value = 4343.3125
$ws.Range("I122").Value = 2113.5
$ws.Range("J122").Value = 6573.125
$ws.Range("K122").Value = 6340.5
$ws.Range("L122").Value = 19719.375
$ws.Range("M122").Value = -3890.5
$ws.Range("N122").Value = -24619.375
$ws.Range("H132").Value = 4730.825
$ws.Range("I132").Value = 2118.5
$ws.Range("J132").Value = 7343.15
$ws.Range("K132").Value = 6355.5
$ws.Range("L132").Value = 22029.45
$ws.Range("M132").Value = -3825.5
$ws.Range("N132").Value = -27089.45
$ws.Range("H134").Value = 35858.062
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 35858.062
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 35858.062
$ws.Range("N134").Value = -45998.062
$ws.Range("H138").Value = 92969.75
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 92969.75
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 92969.75
$ws.Range("N138").Value = -103249.75
$ws.Range("H139").Value = 37226.668
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 37226.668
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 37226.668
$ws.Range("N139").Value = -47506.668
$ws.Range("H140").Value = 92842.86
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 92842.86
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 92842.86
$ws.Range("N140").Value = -103202.86
$ws.Range("H141").Value = 40216.844
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 40216.844
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 40216.844
$ws.Range("N141").Value = -50576.844

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6537260.5
$ws.Range("I132").Value = 838.7111
$ws.Range("J132").Value = 55560424
$ws.Range("K132").Value = 2516.1333
$ws.Range("L132").Value = 166681272
$ws.Range("M132").Value = 13.86670000000004
$ws.Range("N132").Value = -166686332
$ws.Range("H133").Value = 35229.23
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 35229.23
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 35229.23
$ws.Range("N133").Value = -45349.23
$ws.Range("H136").Value = 2508.1667
$ws.Range("I136").Value = 688.6667
$ws.Range("J136").Value = 7966.6665
$ws.Range("K136").Value = 2066.0001
$ws.Range("L136").Value = 23899.9995
$ws.Range("M136").Value = 483.9998999999998
$ws.Range("N136").Value = -28999.9995
$ws.Range("H138").Value = 41499.668
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 41499.668
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 41499.668
$ws.Range("N138").Value = -51779.668
$ws.Range("H139").Value = 35518.95
$ws.Range("I139").Value = 40650
$ws.Range("J139").Value = 35383.92
$ws.Range("K139").Value = 40650
$ws.Range("L139").Value = 35383.92
$ws.Range("M139").Value = -35510
$ws.Range("N139").Value = -45663.92
$ws.Range("H140").Value = 41152.7
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 41152.7
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 41152.7
$ws.Range("N140").Value = -51512.7
$ws.Range("H141").Value = 42285.2
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 42285.2
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 42285.2
$ws.Range("N141").Value = -52645.2
